$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new September 2024 transect data row (row 32)
$ws.Range("A32").Value = 2024
$ws.Range("B32").Value = 9
$ws.Range("C32").Value = 16
$ws.Range("D32").Value = 1
$ws.Range("E32").Value = 0
$ws.Range("F32").Value = 1
$ws.Range("G32").Value = 0
$ws.Range("H32").Value = "presence_abscence"

# Update selection to reflect the next empty row, as in the saved workbook
$ws.Range("A33").Select()
